# Fruta / hortaliza, semanal
# Insert a new weekly price-report block (2 rows) for variety "Andross" at the
# top of the Durazno / Feria Lagunitas de Puerto Montt data range, pushing the
# existing rows (formerly 169-196) down to 171-198.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 169-170 (shifts old 169:196 down to 171:198,
# copying formatting - e.g. the date style on column D - from the row below).
$ws.Rows("169:170").Insert()

# Row 169: Andross / Especial
$ws.Range("A169").Value = 4
$ws.Range("B169").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C169").Value = "Los Lagos"
$ws.Range("D169").Value = 44617
$ws.Range("E169").Value = 10
$ws.Range("F169").Value = "Fruta"
$ws.Range("G169").Value = 100103
$ws.Range("H169").Value = "Frutos de hueso (carozo)"
$ws.Range("I169").Value = 100103004
$ws.Range("J169").Value = "Durazno"
$ws.Range("K169").Value = "Andross"
$ws.Range("L169").Value = "Especial"
$ws.Range("M169").Value = 200
$ws.Range("N169").Value = 19000
$ws.Range("O169").Value = 19000
$ws.Range("P169").Value = 19000
$ws.Range("Q169").Value = "$/caja 15 kilos empedrada"
$ws.Range("R169").Value = "Región de O'Higgins"
$ws.Range("S169").Value = 1267
$ws.Range("T169").Value = 15

# Row 170: Andross / Primera
$ws.Range("A170").Value = 4
$ws.Range("B170").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C170").Value = "Los Lagos"
$ws.Range("D170").Value = 44617
$ws.Range("E170").Value = 10
$ws.Range("F170").Value = "Fruta"
$ws.Range("G170").Value = 100103
$ws.Range("H170").Value = "Frutos de hueso (carozo)"
$ws.Range("I170").Value = 100103004
$ws.Range("J170").Value = "Durazno"
$ws.Range("K170").Value = "Andross"
$ws.Range("L170").Value = "Primera"
$ws.Range("M170").Value = 500
$ws.Range("N170").Value = 16000
$ws.Range("O170").Value = 17000
$ws.Range("P170").Value = 16500
$ws.Range("Q170").Value = "$/caja 15 kilos empedrada"
$ws.Range("R170").Value = "Región de O'Higgins"
$ws.Range("S170").Value = 1100
$ws.Range("T170").Value = 15
